$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New working set of sequence words (column B = "word"), replacing the previous list.
$words = @(
    "sichern",
    "pflanzen",
    "heben",
    "bluten",
    "kratzen",
    "passen",
    "spinnen",
    "weichen",
    "formen",
    "drehen",
    "liegen",
    "fahren",
    "tauchen",
    "trennen",
    "boxen",
    "mögen",
    "starten",
    "lösen",
    "herrschen",
    "sinken",
    "zögern",
    "bellen",
    "wüten",
    "ehren",
    "werfen",
    "wählen",
    "geben",
    "stimmen",
    "suchen",
    "wundern",
    "rauchen",
    "schätzen",
    "wellen",
    "heulen",
    "schenken",
    "knarren",
    "erben",
    "schlagen",
    "kämpfen",
    "platzen",
    "bergen",
    "dienen",
    "kranken",
    "zielen",
    "lassen",
    "kehren",
    "fallen",
    "irren",
    "stärken",
    "stellen",
    "schwingen",
    "stehlen",
    "altern",
    "lenken",
    "scheinen",
    "lügen",
    "runden",
    "schwören",
    "backen",
    "wachsen",
    "proben",
    "feiern",
    "streichen",
    "pfeifen",
    "lockern",
    "stammen",
    "kriegen",
    "tropfen",
    "arten",
    "sprengen",
    "schulden",
    "seufzen",
    "äußern",
    "füttern",
    "dringen",
    "schlucken",
    "fehlen",
    "kichern",
    "albern",
    "spielen",
    "werden",
    "fällen",
    "bieten",
    "warnen",
    "grüßen",
    "gelten",
    "heißen",
    "husten",
    "zünden",
    "posten",
    "schmecken",
    "schultern",
    "filmen",
    "graben",
    "jubeln",
    "achten",
    "liefern",
    "sparen",
    "heilen",
    "spenden",
    "flüchten",
    "führen",
    "schreiten",
    "enden",
    "fragen",
    "jagen",
    "fischen",
    "klappen",
    "folgen",
    "mauern",
    "schrecken",
    "stecken",
    "reizen",
    "sammeln",
    "warten",
    "messen",
    "locken",
    "loben",
    "decken",
    "schließen",
    "ahnen",
    "wenden",
    "klingen",
    "machen",
    "räumen",
    "fordern",
    "schleppen",
    "schämen",
    "sperren",
    "tollen",
    "leeren",
    "merken",
    "saufen",
    "stechen",
    "kümmern",
    "lächeln",
    "haben",
    "zeigen",
    "lohnen",
    "fließen",
    "öffnen",
    "leugnen",
    "helfen",
    "wagen",
    "bitten",
    "quälen",
    "wehtun",
    "bauen",
    "spüren",
    "fühlen",
    "lesen",
    "dürfen",
    "reisen",
    "wollen",
    "streifen",
    "malen",
    "siegen",
    "beißen",
    "biegen",
    "brauchen",
    "töten",
    "löschen",
    "gründen",
    "planen",
    "fangen",
    "rasen",
    "scheitern",
    "sorgen",
    "rechnen",
    "münzen",
    "schwächen",
    "ziehen",
    "betteln",
    "fesseln",
    "ändern",
    "trauen",
    "flehen",
    "hauen",
    "zahlen",
    "kosten",
    "freuen",
    "ärgern",
    "greifen",
    "sterben",
    "nennen",
    "treiben",
    "klettern",
    "rufen",
    "feuern",
    "spannen",
    "wirken"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}

# The previous word list had one more entry than the new one, so the last
# (now unused) data row is removed to keep one row per word.
$lastWordRow = $words.Length + 2
$oldLastRow = 193
if ($oldLastRow -gt ($lastWordRow - 1)) {
    $ws.Range($ws.Cells.Item($lastWordRow, 1), $ws.Cells.Item($oldLastRow, 1)).EntireRow.Delete()
}

Write-Host "Updated word column with new working set of sequences ($($words.Length) words)."
